$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "68.695.60"; E = "  +1.91%  " },
    @{ Row = 3; D = "3.772.45"; E = "  +0.45%  " },
    @{ Row = 4; D = $null; E = "  +0.01%  " },
    @{ Row = 5; D = "597.47"; E = "  +0.47%  " },
    @{ Row = 6; D = "168.71"; E = "  -0.79%  " },
    @{ Row = 7; D = "3.771.87"; E = "  +0.48%  " },
    @{ Row = 8; D = $null; E = "  -0.07%  " },
    @{ Row = 9; D = $null; E = "  -0.34%  " },
    @{ Row = 10; D = $null; E = "  -2.43%  " },
    @{ Row = 11; D = $null; E = "  +0.05%  " },
    @{ Row = 12; D = $null; E = "  -0.51%  " },
    @{ Row = 13; D = $null; E = "  -5.17%  " },
    @{ Row = 14; D = "36.60"; E = "  +0.19%  " },
    @{ Row = 15; D = "4.408.06"; E = "  +0.44%  " },
    @{ Row = 16; D = "3.771.41"; E = "  -0.12%  " },
    @{ Row = 17; D = "68.695.58"; E = "  +1.88%  " },
    @{ Row = 18; D = "18.07"; E = "  -2.74%  " },
    @{ Row = 19; D = "7.05"; E = "  -1.96%  " },
    @{ Row = 20; D = $null; E = "  -0.37%  " },
    @{ Row = 21; D = "10.90"; E = "  +3.60%  " },
    @{ Row = 22; D = "469.66"; E = $null },
    @{ Row = 23; D = "0.705"; E = "  -1.91%  " },
    @{ Row = 24; D = "84.53"; E = $null },
    @{ Row = 25; D = $null; E = "  -0.16%  " },
    @{ Row = 26; D = $null; E = "  +0.36%  " },
    @{ Row = 27; D = "12.14"; E = "  +0.04%  " },
    @{ Row = 28; D = "10.20"; E = "  -1.18%  " },
    @{ Row = 29; D = $null; E = "  +0.05%  " },
    @{ Row = 30; D = "3.920.10"; E = "  +0.28%  " },
    @{ Row = 31; D = $null; E = "  -3.13%  " },
    @{ Row = 32; D = "7.41"; E = "  -3.16%  " },
    @{ Row = 33; D = "30.16"; E = "  -1.26%  " },
    @{ Row = 34; D = $null; E = "  -0.66%  " },
    @{ Row = 35; D = "9.32"; E = "  +2.19%  " },
    @{ Row = 37; D = "3.728.75"; E = "  +0.24%  " },
    @{ Row = 38; D = $null; E = "  -1.94%  " },
    @{ Row = 39; D = "3.48"; E = "  -8.85%  " },
    @{ Row = 40; D = $null; E = "  +1.34%  " },
    @{ Row = 41; D = "1.01"; E = "  +0.80%  " },
    @{ Row = 42; D = $null; E = "  +0.14%  " },
    @{ Row = 43; D = $null; E = "  -0.06%  " },
    @{ Row = 44; D = $null; E = "  +0.00%  " },
    @{ Row = 45; D = $null; E = "  -1.16%  " },
    @{ Row = 46; D = $null; E = "  +1.27%  " },
    @{ Row = 47; D = $null; E = "  +12.23%  " },
    @{ Row = 48; D = $null; E = "  -1.19%  " },
    @{ Row = 49; D = "46.03"; E = "  +0.25%  " },
    @{ Row = 50; D = "397.73"; E = "  +0.04%  " },
    @{ Row = 51; D = "146.06"; E = "  +5.74%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.NumberFormat = "General"
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
